$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated column F values (model-predicted figures) per row, reflecting
# refreshed Stan model output. Only rows with a nonzero F value change;
# rows that were 0 remain 0.

$ws.Range("F4").Value = 10105.66313977776
$ws.Range("F5").Value = 131744.5991485662
$ws.Range("F6").Value = 291323.4506903531
$ws.Range("F7").Value = 444487.4534715106
$ws.Range("F8").Value = 596543.1838339053
$ws.Range("F9").Value = 744605.9023646732
$ws.Range("F10").Value = 881588.6704431179
$ws.Range("F11").Value = 1009967.858438044
$ws.Range("F13").Value = 13568.27480762583
$ws.Range("F14").Value = 60224.0257937783
$ws.Range("F15").Value = 103141.2046088116
$ws.Range("F16").Value = 140905.4605896143
$ws.Range("F17").Value = 170951.4727257513
$ws.Range("F18").Value = 195779.9721990683
$ws.Range("F19").Value = 219407.5135204255
$ws.Range("F20").Value = 242211.932361127
$ws.Range("F21").Value = 262311.784838131
$ws.Range("F43").Value = 28578.66628413394
$ws.Range("F44").Value = 99153.80806321453
$ws.Range("F45").Value = 168933.0664041296
$ws.Range("F46").Value = 232284.3731039814
$ws.Range("F47").Value = 290364.3514974142
$ws.Range("F48").Value = 342749.152422583
$ws.Range("F49").Value = 389735.8485348797
$ws.Range("F50").Value = 432154.1192701443
$ws.Range("F51").Value = 470631.7190954499
$ws.Range("F53").Value = 34620.94955431801
$ws.Range("F54").Value = 134690.1954895044
$ws.Range("F55").Value = 234302.9564667891
$ws.Range("F56").Value = 319565.848698434
$ws.Range("F57").Value = 396100.3448752524
$ws.Range("F58").Value = 462241.2633435281
$ws.Range("F59").Value = 519958.7790941996
$ws.Range("F60").Value = 570839.2983861894
$ws.Range("F61").Value = 617271.110438143
$ws.Range("F64").Value = 2138.565807342776
$ws.Range("F65").Value = 11370.7435639316
$ws.Range("F66").Value = 23541.76004056824
$ws.Range("F67").Value = 37549.02484970899
$ws.Range("F68").Value = 53119.35748416685
$ws.Range("F69").Value = 68865.36846666313
$ws.Range("F70").Value = 84495.33829138298
$ws.Range("F71").Value = 100644.1891555921
$ws.Range("F73").Value = 140815.5151530473
$ws.Range("F74").Value = 573820.6623800445
$ws.Range("F75").Value = 969066.1330866341
$ws.Range("F76").Value = 1299268.821291199
$ws.Range("F77").Value = 1566823.654840165
$ws.Range("F78").Value = 1788833.940871002
$ws.Range("F79").Value = 1980034.902399502
$ws.Range("F80").Value = 2141884.404496072
$ws.Range("F81").Value = 2279532.12610808
$ws.Range("F83").Value = 7739.120050480278
$ws.Range("F84").Value = 75331.32778548842
$ws.Range("F85").Value = 134189.2037464671
$ws.Range("F86").Value = 204837.9871287126
$ws.Range("F87").Value = 263740.1813717776
$ws.Range("F88").Value = 313817.3036328926
$ws.Range("F89").Value = 358502.959583152
$ws.Range("F90").Value = 397273.4995249931
$ws.Range("F91").Value = 431290.3298660588
$ws.Range("F103").Value = 84981.52761576534
$ws.Range("F104").Value = 207071.7256241434
$ws.Range("F105").Value = 306031.0427674481
$ws.Range("F106").Value = 385849.7162298238
$ws.Range("F107").Value = 451733.3239710296
$ws.Range("F108").Value = 506104.4826058967
$ws.Range("F109").Value = 551406.015424223
$ws.Range("F110").Value = 589479.5843029298
$ws.Range("F111").Value = 622303.2462811209
$ws.Range("F115").Value = 899.5048846092382
$ws.Range("F116").Value = 2479.605755001405
$ws.Range("F117").Value = 4108.011092957257
$ws.Range("F118").Value = 6058.172388820716
$ws.Range("F119").Value = 8321.077533350475
$ws.Range("F120").Value = 10593.71796108524
$ws.Range("F121").Value = 12850.156846009
$ws.Range("F123").Value = 44587.16188386771
$ws.Range("F124").Value = 149581.6732378375
$ws.Range("F125").Value = 246934.5975668693
$ws.Range("F126").Value = 332587.6484833481
$ws.Range("F127").Value = 406463.4688868849
$ws.Range("F128").Value = 471082.540776048
$ws.Range("F129").Value = 527163.0013529642
$ws.Range("F130").Value = 576292.5198190317
$ws.Range("F131").Value = 619895.7154822211
$ws.Range("F133").Value = 520.4979779982244
$ws.Range("F134").Value = 10606.78200088063
$ws.Range("F135").Value = 25691.55865322972
$ws.Range("F136").Value = 44014.41636152194
$ws.Range("F137").Value = 62572.68091317185
$ws.Range("F138").Value = 81751.44444225708
$ws.Range("F139").Value = 101335.3062043179
$ws.Range("F140").Value = 120573.3721424613
$ws.Range("F141").Value = 139136.6031172179
$ws.Range("F143").Value = 96262.61298634925
$ws.Range("F144").Value = 415341.0994695302
$ws.Range("F145").Value = 708505.5654050233
$ws.Range("F146").Value = 959159.8047943958
$ws.Range("F147").Value = 1164623.24685351
$ws.Range("F148").Value = 1337532.838319361
$ws.Range("F149").Value = 1484390.51574899
$ws.Range("F150").Value = 1612696.264888573
$ws.Range("F151").Value = 1722952.083070286
$ws.Range("F153").Value = 101314.8331732719
$ws.Range("F154").Value = 260011.0903092448
$ws.Range("F155").Value = 381241.9367681512
$ws.Range("F156").Value = 474182.8495331422
$ws.Range("F157").Value = 545986.0260957096
$ws.Range("F158").Value = 602121.2560361088
$ws.Range("F159").Value = 647354.6836661688
$ws.Range("F160").Value = 685730.35495316
$ws.Range("F161").Value = 717607.7074539758
$ws.Range("F163").Value = 12960.35255315226
$ws.Range("F164").Value = 70867.30265577037
$ws.Range("F165").Value = 149950.2616852272
$ws.Range("F166").Value = 234978.4919618971
$ws.Range("F167").Value = 320173.7718526365
$ws.Range("F168").Value = 401908.3838226532
$ws.Range("F169").Value = 480101.3293867035
$ws.Range("F170").Value = 554029.4837403839
$ws.Range("F171").Value = 624094.8509311627
$ws.Range("F173").Value = 93149.20169396907
$ws.Range("F174").Value = 242618.257618423
$ws.Range("F175").Value = 368602.4665311126
$ws.Range("F176").Value = 472988.9501439421
$ws.Range("F177").Value = 559847.9423868263
$ws.Range("F178").Value = 633581.454738825
$ws.Range("F179").Value = 697307.6129650313
$ws.Range("F180").Value = 752100.7776290084
$ws.Range("F181").Value = 799546.3841404328
$ws.Range("F183").Value = 27034.55961865546
$ws.Range("F184").Value = 88611.6338173751
$ws.Range("F185").Value = 147634.7709000691
$ws.Range("F186").Value = 201823.0122720277
$ws.Range("F187").Value = 250487.3833808773
$ws.Range("F188").Value = 294202.7351436985
$ws.Range("F189").Value = 334266.531291351
$ws.Range("F190").Value = 370535.4984290237
$ws.Range("F191").Value = 403256.9070899478
$ws.Range("F193").Value = 16976.44218567253
$ws.Range("F194").Value = 98775.61511680624
$ws.Range("F195").Value = 185773.2656574085
$ws.Range("F196").Value = 261901.8185475472
$ws.Range("F197").Value = 331688.198993951
$ws.Range("F198").Value = 394072.898115203
$ws.Range("F199").Value = 449314.0938903285
$ws.Range("F200").Value = 498812.8339206895
$ws.Range("F201").Value = 545286.8747045462
$ws.Range("F203").Value = 13194.85177069857
$ws.Range("F204").Value = 77800.79346152839
$ws.Range("F205").Value = 162632.3700006855
$ws.Range("F206").Value = 248924.4526125579
$ws.Range("F207").Value = 334017.7673928849
$ws.Range("F208").Value = 414822.88408156
$ws.Range("F209").Value = 490417.4163615415
$ws.Range("F210").Value = 560862.3870329701
$ws.Range("F211").Value = 626410.4618427065
